$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells (revised AgTests / AgPosit figures) ---
$ws.Range("F354").Value = 316624
$ws.Range("G354").Value = 2889
$ws.Range("F424").Value = 266830
$ws.Range("F431").Value = 171483
$ws.Range("F436").Value = 145560
$ws.Range("F437").Value = 167502
$ws.Range("F445").Value = 84541
$ws.Range("F446").Value = 86967
$ws.Range("F449").Value = 60204
$ws.Range("F451").Value = 86676
$ws.Range("F452").Value = 74598
$ws.Range("F459").Value = 59847
$ws.Range("F465").Value = 61520
$ws.Range("F466").Value = 51057
$ws.Range("F478").Value = 55138
$ws.Range("F498").Value = 9248
$ws.Range("F518").Value = 7273
$ws.Range("F519").Value = 8075
$ws.Range("F520").Value = 10510
$ws.Range("F521").Value = 6956
$ws.Range("F522").Value = 5233
$ws.Range("F523").Value = 10334
$ws.Range("F524").Value = 7901
$ws.Range("F525").Value = 7717
$ws.Range("F526").Value = 8911
$ws.Range("F527").Value = 11678
$ws.Range("F528").Value = 8165
$ws.Range("F529").Value = 5784
$ws.Range("F530").Value = 12967
$ws.Range("G530").Value = 45
$ws.Range("F531").Value = 9335
$ws.Range("F532").Value = 10362
$ws.Range("F533").Value = 11919
$ws.Range("F534").Value = 16848
$ws.Range("F535").Value = 10180
$ws.Range("F536").Value = 8010
$ws.Range("F537").Value = 14035
$ws.Range("F538").Value = 11298
$ws.Range("F539").Value = 10660
$ws.Range("F540").Value = 12507
$ws.Range("F541").Value = 16618
$ws.Range("F542").Value = 10356
$ws.Range("F543").Value = 4705
$ws.Range("F544").Value = 14392
$ws.Range("F545").Value = 16678
$ws.Range("F546").Value = 3901
$ws.Range("F547").Value = 14136
$ws.Range("G547").Value = 152
$ws.Range("F548").Value = 17170
$ws.Range("F549").Value = 10708
$ws.Range("F550").Value = 8513
$ws.Range("G550").Value = 87
$ws.Range("F551").Value = 17824
$ws.Range("F552").Value = 15595
$ws.Range("F553").Value = 15416
$ws.Range("G553").Value = 178
$ws.Range("F554").Value = 17309
$ws.Range("G554").Value = 182
$ws.Range("F555").Value = 21630
$ws.Range("G555").Value = 181
$ws.Range("F556").Value = 12140
$ws.Range("F557").Value = 10897
$ws.Range("G557").Value = 148
$ws.Range("F558").Value = 24942
$ws.Range("G558").Value = 292
$ws.Range("F559").Value = 22485
$ws.Range("G559").Value = 271
$ws.Range("F560").Value = 5970
$ws.Range("G560").Value = 92
$ws.Range("F561").Value = 23813
$ws.Range("G561").Value = 388
$ws.Range("F562").Value = 26724
$ws.Range("G562").Value = 272
$ws.Range("F563").Value = 13886
$ws.Range("G563").Value = 169
$ws.Range("F564").Value = 14053
$ws.Range("G564").Value = 194
$ws.Range("F565").Value = 28443
$ws.Range("G565").Value = 365
$ws.Range("F566").Value = 25852
$ws.Range("G566").Value = 332
$ws.Range("F567").Value = 23081
$ws.Range("G567").Value = 315

# --- Append new daily rows 568-571 ---
$newRows = @(
    @(44462, 406760, 8010, 829, 12594, 23006, 283),
    @(44463, 407799, 10449, 1039, 12596, 28935, 318),
    @(44464, 408488, 6584, 689, 12596, 13203, 198),
    @(44465, 408609, 1518, 121, 12596, 10191, 194)
)

$r = 568
foreach ($rowData in $newRows) {
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $ws.Cells.Item($r, 7).Value = $rowData[6]
    $r++
}
